$d = $word.ActiveDocument

# The document opens with a run reading "EUFEMIO PEREZ" (bold, 12pt,
# minor-theme fonts). We need to insert a brand-new run reading "Test "
# immediately before it, in a separate <w:r> with the same run
# formatting (bold, 12pt / sz 24).

$insertText = "Test "

# Grab the very first character of the document (start of "EUFEMIO PEREZ")
# and insert the new text right before it. Word merges this into the
# existing run since the inherited formatting already matches.
$target = $d.Paragraphs(1).Range.Characters(1)
$target.InsertBefore($insertText)

# Re-select just the text we inserted (now at the very start of the
# document) and nudge its formatting: first to a different size, which
# forces Word to split it back out into its own run, then to the
# correct size/weight as a second, separate operation so the run stays
# split instead of re-merging with the following "EUFEMIO PEREZ" run.
$len = $insertText.Length

$newRun = $d.Range(0, $len)
$newRun.Font.Bold = $true
$newRun.Font.Size = 13

$newRun2 = $d.Range(0, $len)
$newRun2.Font.Bold = $true
$newRun2.Font.Size = 12
